# "Generate Report for Handback" - update the localization-status workbook
# to reflect that the zh-cn and de-de handback packages have been
# generated/returned: Status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns get populated (with a
# hyperlink on the target-file cell, matching the existing source-file
# hyperlink), and a few columns are widened to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/88d0b8bb41488b37b1d102a196874aa80cb06113/e2e/"

$file1 = "bd4c2053-d45e-480a-aca8-056d877ae46d"
$file2 = "e8c2383c-8f44-4b42-8085-804ad6e80005"

# ---------------------------------------------------------------------
# Overview sheet: Status columns for both locales (E/F) move to handed
# back, and those columns need to be wide enough for the longer text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.17
$wsOverview.Columns.Item(6).ColumnWidth = 29.17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Row 2 (bd4c2053...): target + handback file name, handback datetime
$wsZh.Range("I2").Value = "$file1.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$baseUrl$file1.md", "", "", "$file1.md") | Out-Null
$wsZh.Range("J2").Value = "$file1.5f7e76d22d1ae8563ba28a639bb9c5f37d651248.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-29 23:04:53"

# Row 3 (e8c2383c...): target + handback file name, handback datetime
$wsZh.Range("I3").Value = "$file2.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$baseUrl$file2.md", "", "", "$file2.md") | Out-Null
$wsZh.Range("J3").Value = "$file2.7c282366b5d539f6b88b14af72e916d04fdb65b8.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-29 23:04:53"

$wsZh.Columns.Item(3).ColumnWidth = 29.17
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Row 2 (bd4c2053...): target + handback file name, handback datetime
$wsDe.Range("I2").Value = "$file1.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$baseUrl$file1.md", "", "", "$file1.md") | Out-Null
$wsDe.Range("J2").Value = "$file1.5f7e76d22d1ae8563ba28a639bb9c5f37d651248.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-29 23:05:00"

# Row 3 (e8c2383c...): target + handback file name, handback datetime
$wsDe.Range("I3").Value = "$file2.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$baseUrl$file2.md", "", "", "$file2.md") | Out-Null
$wsDe.Range("J3").Value = "$file2.7c282366b5d539f6b88b14af72e916d04fdb65b8.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-29 23:05:00"

$wsDe.Columns.Item(3).ColumnWidth = 29.17
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
